$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.607.42"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "2.934.56"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "198.01"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").Value = "593.54"
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("D10").Value = "2.937.32"
$ws.Range("E10").Value = "  +1.86%  "
$ws.Range("D11").Value = "0.438"
$ws.Range("E11").Value = "  +9.87%  "
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").Value = "3.477.40"
$ws.Range("E13").Value = "  +2.78%  "
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "76.530.16"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "28.29"
$ws.Range("E16").Value = "  +2.76%  "
$ws.Range("D17").Value = "0.0000187"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").Value = "2.955.63"
$ws.Range("E18").Value = "  +3.64%  "
$ws.Range("D19").Value = "13.44"
$ws.Range("E19").Value = "  +6.63%  "
$ws.Range("D20").Value = "8.65"
$ws.Range("E20").Value = "  -3.91%  "
$ws.Range("D21").Value = "372.74"
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("D22").Value = "4.28"
$ws.Range("E22").Value = "  +3.22%  "
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("D24").Value = "72.16"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("D26").Value = "3.094.31"
$ws.Range("E26").Value = "  +2.03%  "
$ws.Range("D27").Value = "4.24"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "9.55"
$ws.Range("E29").Value = "  +0.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").Value = "8.24"
$ws.Range("E31").Value = "  +5.87%  "
$ws.Range("E32").Value = "  -3.49%  "
$ws.Range("D33").Value = "496.06"
$ws.Range("E33").Value = "  -3.23%  "
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "165.28"
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("D37").Value = "0.111"
$ws.Range("E37").Value = "  +19.60%  "
$ws.Range("D38").Value = "20.09"
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.390"
$ws.Range("E39").Value = "  +12.23%  "
$ws.Range("E40").Value = "  +2.21%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.110"
$ws.Range("E41").Value = "  -6.11%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "179.18"
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("E44").Value = "  -3.74%  "
$ws.Range("E45").Value = "  -3.20%  "
$ws.Range("D46").Value = "40.13"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").Value = "1.17"
$ws.Range("E47").Value = "  -5.06%  "
$ws.Range("D48").Value = "0.587"
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("D49").Value = "3.88"
$ws.Range("E49").Value = "  +2.82%  "
$ws.Range("D50").Value = "2.27"
$ws.Range("E50").Value = "  -3.78%  "
$ws.Range("E51").Value = "  -6.06%  "
